$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder header row values (B1:F1), swapping pairs:
#   B1: kitchens_1      -> living_rooms_1
#   C1: living_rooms_1  -> kitchens_1
#   D1: bedrooms_2      -> living_rooms_2
#   E1: kitchens_2      -> bedrooms_2
#   F1: living_rooms_2  -> kitchens_2
$ws.Range("B1").Value = "living_rooms_1"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("D1").Value = "living_rooms_2"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "kitchens_2"
